# Add a detected license column to the export report
#
# OssList.xlsx header row gains a new "Detected License" column, and the
# existing "License" column is relabeled "Declared License".
#
# Layout before:  ... E=OSS Type | F=License | G=License Type | H=Obligation | I=Home Page ...
# Layout after:   ... E=OSS Type | F=Declared License | G=License Type | H=Obligation | I=Detected License | J=Home Page ...

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column in front of the old "Home Page" column (column I);
# this shifts Home Page..Vulnerability right by one and carries the
# surrounding row/column formatting along with it.
$ws.Columns.Item(9).Insert() | Out-Null

# Relabel the existing "License" header as "Declared License".
$ws.Range("F1").Value = "Declared License"

# Label the newly inserted column as "Detected License" (leading space
# matches the source data's header text).
$ws.Range("I1").Value = " Detected License"

# Give the new column the same display width as its left neighbour
# (both were ~15.875 characters wide in the source report).
$ws.Columns.Item(9).ColumnWidth = $ws.Columns.Item(8).ColumnWidth

# Move the active selection the way the author's session left it.
$ws.Range("F10").Select() | Out-Null

Write-Host "Inserted Detected License column and relabeled License -> Declared License"
